$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record (31-Jan-2023) is inserted at row 193, pushing the
# existing data block (rows 193:311) down by one row (193:311 -> 194:312).
$ws.Rows("193:193").Insert()

# The freshly inserted row starts out blank; seed it by duplicating the row
# that is now directly below it (row 195, formerly row 194) and then
# overwrite just the date for the new record.
$ws.Range("A195:T195").Copy($ws.Range("A193:T193"))
$ws.Range("D193").Value = 44957
